# TC41_Canine_Filter_Breed-Samoyed.xlsx - "startup" sheet
# Fix the CasesTab (row 2) Cypher query in column B: it was erroneously
# returning a `Cohort` column (joined via a `(co:cohort)` match that isn't
# part of this query's intent) - drop that trailing RETURN item/newline so
# the Cases query matches the Samples/Files queries' shape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesQueryLines = @(
    "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)",
    "WHERE demo.breed  IN ['Samoyed']",
    "MATCH (c)<--(diag:diagnosis)",
    "OPTIONAL MATCH (samp:sample)-->(c)",
    "OPTIONAL MATCH (co:cohort)<-[*]-(c)",
    "WITH DISTINCT c, s, demo, diag, co",
    "RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,",
    "        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,",
    "        coalesce(s.clinical_study_type, '') AS  ``Study Type``,",
    "        coalesce(demo.breed, '') AS Breed ,",
    "        coalesce(diag.disease_term, '') AS Diagnosis ,",
    "        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,",
    "        coalesce(demo.patient_age_at_enrollment, '') AS Age ,",
    "        coalesce(demo.sex, '') AS Sex ,",
    "        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,",
    "        coalesce(demo.weight, '') AS ``Weight (kg)``,",
    "        coalesce(diag.best_response, '') AS ``Response to Treatment``"
)
$newCasesQuery = [string]::Join("`r`n", $casesQueryLines)

$ws.Range("B2").Value = $newCasesQuery

# Match the author's final active selection on this sheet (B2).
$ws.Range("B2").Select() | Out-Null
